# Updated cryptos list on Sun Mar 31 03:55:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-CellText $ws "D2" "69.733.31"
Set-CellText $ws "E2" "  +0.01%  "
Set-CellText $ws "D3" "3.532.06"
Set-CellText $ws "E3" "  +1.23%  "
Set-CellText $ws "E4" "  -0.02%  "
Set-CellText $ws "D5" "606.73"
Set-CellText $ws "E5" "  +0.03%  "
Set-CellText $ws "D6" "194.14"
Set-CellText $ws "E6" "  +0.86%  "
Set-CellText $ws "D7" "0.622"
Set-CellText $ws "E7" "  -0.34%  "
Set-CellText $ws "E8" "  +0.04%  "
Set-CellText $ws "D9" "0.201"
Set-CellText $ws "E9" "  -4.94%  "
Set-CellText $ws "D10" "0.645"
Set-CellText $ws "E10" "  -2.07%  "
Set-CellText $ws "D11" "53.20"
Set-CellText $ws "E11" "  -0.26%  "
Set-CellText $ws "E12" "  -1.31%  "
Set-CellText $ws "D13" "9.45"
Set-CellText $ws "E13" "  -1.37%  "
Set-CellText $ws "D14" "4.092.09"
Set-CellText $ws "E14" "  +0.94%  "
Set-CellText $ws "D15" "592.33"
Set-CellText $ws "E15" "  -2.21%  "
Set-CellText $ws "D16" "12.76"
Set-CellText $ws "E16" "  +0.92%  "
Set-CellText $ws "D17" "69.895.97"
Set-CellText $ws "E17" "  +0.06%  "
Set-CellText $ws "D18" "18.94"
Set-CellText $ws "E18" "  +0.70%  "
Set-CellText $ws "D19" "3.535.43"
Set-CellText $ws "E19" "  +1.54%  "
Set-CellText $ws "E20" "  +1.76%  "
Set-CellText $ws "D21" "0.981"
Set-CellText $ws "E21" "  -0.76%  "
Set-CellText $ws "D22" "17.63"
Set-CellText $ws "E22" "  -0.54%  "
Set-CellText $ws "D23" "102.85"
Set-CellText $ws "E23" "  -2.32%  "
Set-CellText $ws "E24" "  +1.10%  "
Set-CellText $ws "E25" "  -0.12%  "
Set-CellText $ws "E26" "  -0.85%  "
Set-CellText $ws "D27" "10.73"
Set-CellText $ws "E27" "  -1.99%  "
Set-CellText $ws "D28" "9.48"
Set-CellText $ws "E28" "  -3.71%  "
Set-CellText $ws "D29" "33.07"
Set-CellText $ws "D30" "7.01"
Set-CellText $ws "E30" "  -1.75%  "
Set-CellText $ws "D31" "4.21"
Set-CellText $ws "E31" "  -2.31%  "
Set-CellText $ws "D32" "12.27"
Set-CellText $ws "E32" "  -2.98%  "
Set-CellText $ws "E33" "  -0.24%  "
Set-CellText $ws "D34" "63.30"
Set-CellText $ws "E34" "  -1.41%  "
Set-CellText $ws "D35" "3.819.25"
Set-CellText $ws "E35" "  +2.76%  "
Set-CellText $ws "D36" "3.18"
Set-CellText $ws "E36" "  +3.91%  "
Set-CellText $ws "D37" "0.0₃0815"
Set-CellText $ws "E37" "  +2.86%  "
Set-CellText $ws "E38" "  +0.16%  "
Set-CellText $ws "D39" "512.86"
Set-CellText $ws "E39" "  -0.93%  "
Set-CellText $ws "D40" "0.388"
Set-CellText $ws "E40" "  -0.33%  "
Set-CellText $ws "E41" "  -0.81%  "
Set-CellText $ws "D42" "36.42"
Set-CellText $ws "E42" "  -0.30%  "
Set-CellText $ws "E43" "  -2.16%  "
Set-CellText $ws "D44" "0.0446"
Set-CellText $ws "E44" "  -3.32%  "
Set-CellText $ws "D45" "3.37"
Set-CellText $ws "E45" "  +1.64%  "
Set-CellText $ws "E46" "  -1.00%  "
Set-CellText $ws "E47" "  -1.71%  "
Set-CellText $ws "E48" "  +0.13%  "
Set-CellText $ws "D49" "8.46"
Set-CellText $ws "E49" "  -3.09%  "
Set-CellText $ws "B50" "Mantle"
Set-CellText $ws "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText $ws "D50" "1.32"
Set-CellText $ws "E50" "  +1.69%  "
Set-CellText $ws "B51" "FLOKI"
Set-CellText $ws "C51" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-CellText $ws "D51" "0.000244"
Set-CellText $ws "E51" "  +3.07%  "
